$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numSalle (column B) values for rows 2-20
$values = @(30,31,32,33,34,35,36,30,31,32,33,34,35,36,30,31,32,33,34)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the selection on the sheet
$ws.Range("F21").Select()
